$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 4959.7104
$ws.Range("J32").Value = 4900
$ws.Range("L32").Value = 4900
$ws.Range("N32").Value = -5552
# Row 100
$ws.Range("H100").Value = 2673
$ws.Range("I100").Value = 1991.6
$ws.Range("J100").Value = 3240.8333
$ws.Range("K100").Value = 1991.6
$ws.Range("L100").Value = 3240.8333
$ws.Range("M100").Value = -1450.6
$ws.Range("N100").Value = -4322.8333
# Row 113
$ws.Range("H113").Value = 3214.2856
# Row 116
$ws.Range("H116").Value = 6571
$ws.Range("I116").Value = 4249.25
$ws.Range("J116").Value = 9666.666999999999
$ws.Range("K116").Value = 4249.25
$ws.Range("L116").Value = 9666.666999999999
$ws.Range("M116").Value = -807.25
$ws.Range("N116").Value = -16550.667
# Row 125
$ws.Range("H125").Value = 1983
$ws.Range("I125").Value = 1975
$ws.Range("K125").Value = 17775
$ws.Range("M125").Value = -15315
# Row 138
$ws.Range("H138").Value = 4287.477
$ws.Range("J138").Value = 6234.654
$ws.Range("L138").Value = 18703.962
$ws.Range("N138").Value = -28983.962

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4042.7058
$ws.Range("I32").Value = 2578.3403
$ws.Range("K32").Value = 2578.3403
$ws.Range("M32").Value = -2291.3403
# Row 45
$ws.Range("H45").Value = 1304
$ws.Range("I45").Value = 1304
$ws.Range("K45").Value = 1304
$ws.Range("M45").Value = -927
# Row 74
$ws.Range("H74").Value = 1031.4
$ws.Range("I74").Value = 1031.4
$ws.Range("K74").Value = 1031.4
$ws.Range("M74").Value = -157.4000000000001
# Row 77
$ws.Range("H77").Value = 1031.4
$ws.Range("I77").Value = 1031.4
$ws.Range("K77").Value = 5157
$ws.Range("M77").Value = -789
# Row 122
$ws.Range("H122").Value = 2056.818
$ws.Range("I122").Value = 2062.5
$ws.Range("K122").Value = 6187.5
$ws.Range("M122").Value = -3737.5
# Row 132
$ws.Range("H132").Value = 1401.8572
$ws.Range("I132").Value = 1347.8182
$ws.Range("K132").Value = 4043.4546
$ws.Range("M132").Value = -1513.4546

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 111.5
$ws.Range("I7").Value = 53.692307
$ws.Range("K7").Value = 53.692307
$ws.Range("M7").Value = 59.307693
# Row 16
$ws.Range("H16").Value = 4528.8
$ws.Range("I16").Value = 4528.8
$ws.Range("K16").Value = 4528.8
$ws.Range("M16").Value = -4241.8
# Row 22
$ws.Range("H22").Value = 549.3333
$ws.Range("I22").Value = 549.3333
$ws.Range("K22").Value = 549.3333
$ws.Range("M22").Value = -199.3333
# Row 31
$ws.Range("H31").Value = 3232.4666
$ws.Range("I31").Value = 2809.7778
$ws.Range("J31").Value = 7036.6665
$ws.Range("K31").Value = 2809.7778
$ws.Range("L31").Value = 7036.6665
$ws.Range("M31").Value = -2514.7778
$ws.Range("N31").Value = -7626.6665
# Row 34
$ws.Range("H34").Value = 3232.4666
$ws.Range("I34").Value = 2809.7778
$ws.Range("J34").Value = 7036.6665
$ws.Range("K34").Value = 2809.7778
$ws.Range("L34").Value = 7036.6665
$ws.Range("M34").Value = -2607.7778
$ws.Range("N34").Value = -7440.6665
# Row 113
$ws.Range("H113").Value = 4528.8
$ws.Range("I113").Value = 4528.8
$ws.Range("K113").Value = 4528.8
$ws.Range("M113").Value = -2358.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 1451.1111
$ws.Range("J131").Value = 1583.8695
$ws.Range("L131").Value = 4751.6085
$ws.Range("N131").Value = -14831.6085

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 755.5
$ws.Range("I2").Value = 80
$ws.Range("K2").Value = 80
$ws.Range("M2").Value = 33
# Row 3
$ws.Range("H3").Value = 52777.105
$ws.Range("I3").Value = 55698.055
$ws.Range("J3").Value = 200
$ws.Range("K3").Value = 55698.055
$ws.Range("L3").Value = 200
$ws.Range("M3").Value = -55582.055
$ws.Range("N3").Value = -432
# Row 5
$ws.Range("H5").Value = 25
$ws.Range("I5").Value = 25
$ws.Range("K5").Value = 25
$ws.Range("M5").Value = 87
# Row 24
$ws.Range("H24").Value = 246079.77
$ws.Range("I24").Value = 500750
$ws.Range("J24").Value = 10999.538
$ws.Range("K24").Value = 500750
$ws.Range("L24").Value = 10999.538
$ws.Range("M24").Value = -500577
$ws.Range("N24").Value = -11345.538
# Row 41
$ws.Range("H41").Value = 2062.75
$ws.Range("I41").Value = 2062.75
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 2062.75
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -1707.75
$ws.Range("N41").ClearContents()
# Row 92
$ws.Range("H92").Value = 5746.2856
$ws.Range("J92").Value = 5746.2856
$ws.Range("L92").Value = 5746.2856
$ws.Range("N92").Value = -9490.285599999999
# Row 93
$ws.Range("H93").Value = 25768.908
$ws.Range("J93").Value = 25768.908
$ws.Range("L93").Value = 25768.908
$ws.Range("N93").Value = -29512.908
# Row 109
$ws.Range("H109").Value = 6571.4287
$ws.Range("J109").Value = 6571.4287
$ws.Range("L109").Value = 6571.4287
$ws.Range("N109").Value = -8651.4287
# Row 113
$ws.Range("H113").Value = 2196
$ws.Range("I113").Value = 1294
$ws.Range("K113").Value = 1294
$ws.Range("M113").Value = 876
# Row 122
$ws.Range("H122").Value = 62049.766
$ws.Range("I122").Value = 2931
$ws.Range("J122").Value = 128558.375
$ws.Range("K122").Value = 8793
$ws.Range("L122").Value = 385675.125
$ws.Range("M122").Value = -6343
$ws.Range("N122").Value = -390575.125
# Row 123
$ws.Range("H123").Value = 22000.223
$ws.Range("J123").Value = 22000.223
$ws.Range("L123").Value = 22000.223
$ws.Range("N123").Value = -26900.223
# Row 126
$ws.Range("H126").Value = 4407.2
$ws.Range("I126").Value = 3674
$ws.Range("K126").Value = 11022
$ws.Range("M126").Value = -8552

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 3963.2693
$ws.Range("I46").Value = 2277.889
$ws.Range("K46").Value = 2277.889
$ws.Range("M46").Value = -2089.889
# Row 61
$ws.Range("H61").Value = 4922.5835
$ws.Range("I61").Value = 5382.5
$ws.Range("J61").Value = 4002.75
$ws.Range("K61").Value = 5382.5
$ws.Range("L61").Value = 4002.75
$ws.Range("M61").Value = -5180.5
$ws.Range("N61").Value = -4406.75
# Row 113
$ws.Range("H113").Value = 4922.5835
$ws.Range("I113").Value = 5382.5
$ws.Range("J113").Value = 4002.75
$ws.Range("K113").Value = 5382.5
$ws.Range("L113").Value = 4002.75
$ws.Range("M113").Value = -3212.5
$ws.Range("N113").Value = -8342.75
# Row 132
$ws.Range("H132").Value = 2360.2964
$ws.Range("I132").Value = 2201.5454
$ws.Range("K132").Value = 6604.6362
$ws.Range("M132").Value = -4074.6362

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 18
$ws.Range("H18").Value = 8780
$ws.Range("I18").Value = 2900
$ws.Range("J18").Value = 10250
$ws.Range("K18").Value = 2900
$ws.Range("L18").Value = 10250
$ws.Range("M18").Value = -2727
$ws.Range("N18").Value = -10596
# Row 62
$ws.Range("H62").Value = 7617.1177
$ws.Range("I62").Value = 5833.3335
$ws.Range("J62").Value = 7999.357
$ws.Range("K62").Value = 5833.3335
$ws.Range("L62").Value = 7999.357
$ws.Range("M62").Value = -5209.3335
$ws.Range("N62").Value = -9247.357
# Row 65
$ws.Range("H65").Value = 7617.1177
$ws.Range("I65").Value = 5833.3335
$ws.Range("J65").Value = 7999.357
$ws.Range("K65").Value = 29166.6675
$ws.Range("L65").Value = 39996.785
$ws.Range("M65").Value = -26046.6675
$ws.Range("N65").Value = -46236.785
